$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (the numeric index 0..4) is dropped entirely; the former column B
# shifts into A and the former column C shifts into B.
$ws.Range("A1:A5").Delete()

# The only value that doesn't simply carry over from the old column C is the
# last row: former C5 was "7", but the new B5 must read "9". Stage the text
# "9" in a scratch cell via a formula (so it lands as a genuine text value,
# not an auto-coerced number), copy it across as a value, then tidy the
# scratch cell back up.
$ws.Range("Z1").Formula = '="9"'
$ws.Range("Z1").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$excel.CutCopyMode = 0
